# Apply the edits described in the commit diff:
# - A7 ("am-am") -> "am-ha"
# - C7 ("amidou") -> "halidou"
# - Active selection moves to A7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "am-ha"
$ws.Range("C7").Value = "halidou"

$ws.Range("A7").Select()
